$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Account numbers (Conta, column A) whose rows must be removed.
$accountsToRemove = @(
    "000806386",  # Fernanda  50533.82
    "004444380",  # Marcelo   49648.92
    "005599726",  # Jorge     26275.75
    "004342617",  # Juraci     7669.90
    "004948033",  # Guilherme  2260.42
    "005061124"   # Bruno      2204.13
)

foreach ($acct in $accountsToRemove) {
    $found = $ws.Columns("A:A").Find($acct)
    if ($found -ne $null) {
        $ws.Rows($found.Row).Delete()
    }
}
